$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# Table (graphicFrame id=32): mean-number summary rows
$tbl32 = (Get-ShapeById $s 32).Table
$tbl32.Cell(1,1).Shape.TextFrame.TextRange.Text = "Number of paired, normalized concentration and temperature values"
$tbl32.Cell(2,1).Shape.TextFrame.TextRange.Text = "Mean number of paired, normalized concentration and relative humidity values"

# Table (graphicFrame id=75): FRM/FEM goal concentration threshold 25 -> 40
$tbl75 = (Get-ShapeById $s 75).Table
$cell = $tbl75.Cell(2,1)
$tr = $cell.Shape.TextFrame.TextRange
$sub = $tr.Characters(83, 5)
$sub.Text = "≥ 40 "
